$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8015
$ws1.Range("F10").Value = 470
$ws1.Range("F17").Value = 5910
$ws1.Range("F18").Value = 188
$ws1.Range("F19").Value = 276
$ws1.Range("F20").Value = 1918
$ws1.Range("F22").Value = 29
$ws1.Range("F24").Value = 408

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8015
$ws4.Range("F10").Value = 470
$ws4.Range("F18").Value = 5910
$ws4.Range("F20").Value = 188
$ws4.Range("F21").Value = 276
$ws4.Range("F22").Value = 1918
$ws4.Range("F24").Value = 29
$ws4.Range("F26").Value = 408
